$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# --- Rename the two sheets we keep, drop the third ---
$wsNhap = $wb.Worksheets.Item("nhap-linhkien")
$wsNhap.Name = "nhap-thanhpham"

$wsXuat = $wb.Worksheets.Item("xuat-linhkien")
$wsXuat.Name = "xuat-thanhpham"

$wsTon = $wb.Worksheets.Item("ton-linhkien")
$wsTon.Delete()

# --- Shrink the table from 10 columns (A:J) down to 6 (A:F) ---
$wsNhap.Columns("G:J").Delete()
$wsXuat.Columns("G:J").Delete()

# --- New header row for both remaining sheets ---
$headers = @("Tên Hàng", "MCU", "Sổ Hợp Đồng", "Chip", "Ngày Nhập", "Số Lượng")

$wsNhap.Range("A1").Value = $headers[0]
$wsNhap.Range("B1").Value = $headers[1]
$wsNhap.Range("C1").Value = $headers[2]
$wsNhap.Range("D1").Value = $headers[3]
$wsNhap.Range("E1").Value = $headers[4]
$wsNhap.Range("F1").Value = $headers[5]

$wsXuat.Range("A1").Value = $headers[0]
$wsXuat.Range("B1").Value = $headers[1]
$wsXuat.Range("C1").Value = $headers[2]
$wsXuat.Range("D1").Value = $headers[3]
$wsXuat.Range("E1").Value = $headers[4]
$wsXuat.Range("F1").Value = $headers[5]

# --- Sample data row added to "nhap-thanhpham" only ---
$wsNhap.Range("A2").Value = "LED Green"
$wsNhap.Range("B2").Value = "mcu01"
$wsNhap.Range("C2").Value = "sohopdong01"
$wsNhap.Range("D2").Value = "chip01"
$wsNhap.Range("E2").NumberFormat = "@"
$wsNhap.Range("E2").Value = "2021-09-02"
$wsNhap.Range("E2").ClearFormats()
$wsNhap.Range("F2").Value = 12

# --- Restore the original active sheet / selection ---
$wsNhap.Activate()
$wsNhap.Range("A1").Select()
